# Rename the inline picture "Name" metadata (wp:docPr/@name) for the three
# logo images that live in the document's footers/header:
#   - Footer 1 (odd pages)  : Pearson logo  "image2.png" -> "image1.png"
#   - Footer 2 (even pages) : Pearson logo  "image2.png" -> "image1.png"
#   - Header 2 (even pages) : BTEC logo     "image1.jpg" -> "image2.jpg"

$d = $word.ActiveDocument
$sec = $d.Sections.Item(1)

# --- Footer 1 : Pearson Edexcel logo ---------------------------------
$footer1 = $sec.Footers.Item(1)
$pearsonShape1 = $footer1.Range.InlineShapes.Item(1)
$pearsonShape1.Name = "image1.png"

# --- Footer 2 : Pearson Edexcel logo ---------------------------------
$footer2 = $sec.Footers.Item(2)
$pearsonShape2 = $footer2.Range.InlineShapes.Item(1)
$pearsonShape2.Name = "image1.png"

# --- Header 2 : BTEC logo --------------------------------------------
$header2 = $sec.Headers.Item(2)
$btecShape = $header2.Range.InlineShapes.Item(1)
$btecShape.Name = "image2.jpg"

Write-Output "Renamed inline shapes in footer1, footer2 and header2."
